# Budget workbook update: add the "Power" subsection row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 6: label "Power" in column A, value 10 in column B (keeps the
# existing currency-formatted style already present on B6/B7/...).
$ws.Range("A6").Value = "Power"
$ws.Range("B6").Value = 10

# Move/update the active selection to C6, matching the saved workbook state.
$null = $ws.Range("C6").Select()
